# Commit: "Removed Pluto & Sedna from Planets List"
#
# The planets list has one planet name per paragraph. We need to delete the
# trailing "Pluto" and "Sedna" paragraphs entirely. The "Sedna" paragraph
# also happens to hold the document's "_GoBack" bookmark (wrapped in
# w:proofErr spell-check markers); Word's normal behaviour when you delete
# through a bookmark like this is to keep the (now empty) bookmark and let
# it collapse onto the end of whatever text now precedes it -- here, onto
# the end of the "Neptune " paragraph, which becomes the new last paragraph.

$d = $word.ActiveDocument

# Locate the "Pluto" / "Sedna" paragraphs by content instead of a hard-coded
# index, and remember the paragraph immediately before them (the one the
# bookmark should end up attached to).
$targets = @("Pluto", "Sedna")
$firstTargetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs($i).Range.Text.Trim()
    if ($targets -contains $text) {
        if ($firstTargetIndex -eq -1) {
            $firstTargetIndex = $i
        }
    }
}

$precedingIndex = $firstTargetIndex - 1
$lastIndex = $d.Paragraphs.Count

# Delete the target paragraphs from the back forward so earlier indices
# stay valid while we work.
for ($i = $lastIndex; $i -ge $firstTargetIndex; $i--) {
    $d.Paragraphs($i).Range.Delete()
}

# The bookmark belongs right after the preceding paragraph's text (i.e.
# right before the paragraph mark that now ends what is the last paragraph
# in the document).
$preceding = $d.Paragraphs($precedingIndex)
$bookmarkPos = $preceding.Range.End - 1

# Re-adding a bookmark collapsed at the very last character of the whole
# document mis-anchors in this host, so temporarily append a placeholder
# character after it, add the bookmark, then remove the placeholder again.
$endOfDoc = $d.Content.End
$d.Range($endOfDoc, $endOfDoc).InsertAfter("X")

$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))

$d.Range($bookmarkPos, $bookmarkPos + 1).Delete()
